$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.654.34"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.111.23"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  +1.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "348.78"
$ws.Range("E5").Value = "  +3.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.013"
$ws.Range("E6").Value = "  +1.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5262"
$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4507"
$ws.Range("E8").Value = "  -2.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.64"
$ws.Range("E9").Value = "  +2.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08987"
$ws.Range("E10").Value = "  +0.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.169"
$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.40"
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.106.56"
$ws.Range("E13").Value = "  +0.81%  "

$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.013"
$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.29"
$ws.Range("E16").Value = "  +3.05%  "

$ws.Range("E17").Value = "  +4.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.015"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06716"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.33"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("E21").Value = "  +1.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.298"
$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.699.08"
$ws.Range("E23").Value = "  +0.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.83"
$ws.Range("E24").Value = "  +3.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.388"
$ws.Range("E25").Value = "  +0.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.355.38"
$ws.Range("E26").Value = "  +0.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.35"
$ws.Range("E27").Value = "  +0.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.18"
$ws.Range("E28").Value = "  +0.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.526"
$ws.Range("E29").Value = "  -1.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.45"
$ws.Range("E30").Value = "  +2.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.186"
$ws.Range("E31").Value = "  -0.90%  "

$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.632"
$ws.Range("E33").Value = "  -3.70%  "

$ws.Range("E34").Value = "  +3.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.017"
$ws.Range("E35").Value = "  +2.50%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.901"
$ws.Range("E36").Value = "  +6.42%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.20"
$ws.Range("E37").Value = "  -2.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02629"
$ws.Range("E38").Value = "  +2.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06828"
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2305"
$ws.Range("E40").Value = "  +0.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.54"
$ws.Range("E41").Value = "  -2.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6859"
$ws.Range("E42").Value = "  -0.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.281"
$ws.Range("E43").Value = "  +2.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.75"
$ws.Range("E44").Value = "  +5.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.315"
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6405"
$ws.Range("E46").Value = "  +0.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.755"
$ws.Range("E47").Value = "  +2.50%  "

$ws.Range("E48").Value = "  -0.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.243"
$ws.Range("E49").Value = "  -0.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07288"
$ws.Range("E50").Value = "  +2.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.194"
$ws.Range("E51").Value = "  -1.55%  "
